$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.191131666666666
$ws.Range("H2").Value = 3.573395
$ws.Range("I2").Value = 0.02720036629735778
$ws.Range("J2").Value = 0.02720036629735778
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 59.45197733333333
$ws.Range("N2").Value = 178.355932
$ws.Range("O2").Value = 0.304222453049858
$ws.Range("P2").Value = 0.304222453049858
$ws.Range("Q2").Value = 70.81513284768221
$ws.Range("R2").Value = 637.33619562914
$ws.Range("S2").Value = 0.008274962158836866
$ws.Range("T2").Value = 0.008274962158836866

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.191131666666666
$ws.Range("H3").Value = 3.573395
$ws.Range("I3").Value = 0.02720036629735778
$ws.Range("J3").Value = 0.02720036629735778
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 107.1770123333333
$ws.Range("N3").Value = 321.531037
$ws.Range("O3").Value = 0.548436823552382
$ws.Range("P3").Value = 0.5484368235523819
$ws.Range("Q3").Value = 127.6619333289572
$ws.Range("R3").Value = 1148.957399960615
$ws.Range("S3").Value = 0.01491768249158417
$ws.Range("T3").Value = 0.01491768249158416

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.191131666666666
$ws.Range("H4").Value = 3.573395
$ws.Range("I4").Value = 0.02720036629735778
$ws.Range("J4").Value = 0.02720036629735778
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 28.793724
$ws.Range("N4").Value = 86.38117199999999
$ws.Range("O4").Value = 0.1473407233977601
$ws.Range("P4").Value = 0.1473407233977601
$ws.Range("Q4").Value = 34.29711645765999
$ws.Range("R4").Value = 308.6740481189399
$ws.Range("S4").Value = 0.004007721646936748
$ws.Range("T4").Value = 0.004007721646936747

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 34.415161
$ws.Range("H5").Value = 103.245483
$ws.Range("I5").Value = 0.7858954736735307
$ws.Range("J5").Value = 0.7858954736735306
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 59.45197733333333
$ws.Range("N5").Value = 178.355932
$ws.Range("O5").Value = 0.304222453049858
$ws.Range("P5").Value = 0.304222453049858
$ws.Range("Q5").Value = 2046.049371695017
$ws.Range("R5").Value = 18414.44434525516
$ws.Range("S5").Value = 0.2390870488417416
$ws.Range("T5").Value = 0.2390870488417416

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 34.415161
$ws.Range("H6").Value = 103.245483
$ws.Range("I6").Value = 0.7858954736735307
$ws.Range("J6").Value = 0.7858954736735306
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 107.1770123333333
$ws.Range("N6").Value = 321.531037
$ws.Range("O6").Value = 0.548436823552382
$ws.Range("P6").Value = 0.5484368235523819
$ws.Range("Q6").Value = 3688.514134950653
$ws.Range("R6").Value = 33196.62721455587
$ws.Range("S6").Value = 0.4310140172257059
$ws.Range("T6").Value = 0.4310140172257057

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 34.415161
$ws.Range("H7").Value = 103.245483
$ws.Range("I7").Value = 0.7858954736735307
$ws.Range("J7").Value = 0.7858954736735306
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 28.793724
$ws.Range("N7").Value = 86.38117199999999
$ws.Range("O7").Value = 0.1473407233977601
$ws.Range("P7").Value = 0.1473407233977601
$ws.Range("Q7").Value = 990.940647249564
$ws.Range("R7").Value = 8918.465825246076
$ws.Range("S7").Value = 0.1157944076060833
$ws.Range("T7").Value = 0.1157944076060833

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 8.184723
$ws.Range("H8").Value = 24.554169
$ws.Range("I8").Value = 0.1869041600291116
$ws.Range("J8").Value = 0.1869041600291116
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 59.45197733333333
$ws.Range("N8").Value = 178.355932
$ws.Range("O8").Value = 0.304222453049858
$ws.Range("P8").Value = 0.304222453049858
$ws.Range("Q8").Value = 486.597966275612
$ws.Range("R8").Value = 4379.381696480508
$ws.Range("S8").Value = 0.05686044204927954
$ws.Range("T8").Value = 0.05686044204927954

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 8.184723
$ws.Range("H9").Value = 24.554169
$ws.Range("I9").Value = 0.1869041600291116
$ws.Range("J9").Value = 0.1869041600291116
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 107.1770123333333
$ws.Range("N9").Value = 321.531037
$ws.Range("O9").Value = 0.548436823552382
$ws.Range("P9").Value = 0.5484368235523819
$ws.Range("Q9").Value = 877.2141579159169
$ws.Range("R9").Value = 7894.927421243253
$ws.Range("S9").Value = 0.1025051238350921
$ws.Range("T9").Value = 0.102505123835092

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 8.184723
$ws.Range("H10").Value = 24.554169
$ws.Range("I10").Value = 0.1869041600291116
$ws.Range("J10").Value = 0.1869041600291116
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 28.793724
$ws.Range("N10").Value = 86.38117199999999
$ws.Range("O10").Value = 0.1473407233977601
$ws.Range("P10").Value = 0.1473407233977601
$ws.Range("Q10").Value = 235.668655078452
$ws.Range("R10").Value = 2121.017895706068
$ws.Range("S10").Value = 0.02753859414474002
$ws.Range("T10").Value = 0.02753859414474001

